$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows: A,B,C,D,E,F,G,H,I,J
$data = @(
    @(0, 0.03816223951913833, 0.005706788595714274, 0.4603888076369064, 408, 896, 752, 693, 'ALB US Equity_CINR US Equity', $null),
    @(1, 0.1962950046862848, 0.002397393255829661, 0.1696479322062048, 554, 750, 1196, 1213, 'ALB US Equity_CBT US Equity', $null),
    @(2, -0.1102038653647273, -0.001058206357413716, 0.72149012486797, 523, 781, 1239, 1170, $null, 'AI FP Equity_AVD US Equity'),
    @(3, -0.2893718208605331, 0.01226109549120191, 0.5162236589715865, 499, 805, 790, 655, 'AI FP Equity_CINR US Equity', $null),
    @(4, -0.1791224894238869, -0.005022278932140356, 0.06774009136924396, 656, 648, 736, 1673, 'AI FP Equity_HXL US Equity', $null),
    @(5, 0.08833734460097675, 0.002282435514504133, 0.4202462126203527, 581, 723, 1338, 1071, $null, 'ASH US Equity_AVD US Equity'),
    @(6, -0.08614101638564375, 0.006167988539737657, 1.368546857399685, 497, 807, 1225, 1184, $null, 'APD US Equity_AVD US Equity'),
    @(7, 0.01386473180607162, 0.01115909690530437, 0.4869103089271438, 588, 716, 1475, 934, $null, 'AVY US Equity_AVD US Equity'),
    @(8, 0.06614079833792563, 0.006265813580010215, 0.0339659356398244, 612, 692, 1218, 1191, 'AVD US Equity_ADM US Equity', $null),
    @(9, 0.04764228887190058, 0.006562950339515261, 0.05135530613582681, 707, 597, 623, 519, 'AVD US Equity_AXTA US Equity', $null),
    @(10, -0.09830159397202665, 0.007632346204695217, 0.01206840221748883, 640, 664, 1229, 1180, 'AVD US Equity_BAS GR Equity', $null),
    @(11, -0.2288232786991248, 0.009695302544658224, 0.0495051194706272, 583, 721, 782, 663, 'AVD US Equity_CINR US Equity', $null),
    @(12, -0.0168960073888057, -0.009965754992696474, 0.02816644030263519, 630, 674, 438, 802, 'AVD US Equity_IMCD NA Equity', $null),
    @(13, 0.04058718068102718, -0.004665422293643084, 0.03701763287162237, 640, 664, 1001, 1408, 'AVD US Equity_DCI US Equity', $null),
    @(14, 0.04676129608862167, -0.01594678751257761, 0.02224192401974126, 729, 575, 898, 1511, 'AVD US Equity_DSM NA Equity', $null),
    @(15, -0.1262673089474124, 0.00308434225887888, 0.01971033593957043, 629, 675, 1204, 1205, 'AVD US Equity_CBT US Equity', $null),
    @(16, -0.02126142153227906, -0.009701428956883085, 0.01428068362350381, 657, 647, 899, 1510, 'AVD US Equity_CE US Equity', $null),
    @(17, -0.2517129646993059, 0.005972534005781016, 0.01082386483981805, 583, 721, 1158, 1251, 'AVD US Equity_DD US Equity', $null),
    @(18, -0.1202620371574051, 0.006636609013132455, 0.01387790709337839, 93, 112, 6, 0, 'AVD US Equity_DOW US Equity', $null),
    @(19, 0.1108204848441618, -0.004127778065523446, 0.01346231592237077, 632, 672, 1132, 1277, 'AVD US Equity_ECL US Equity', $null),
    @(20, -0.09384222836719225, 0.003281269931679098, 0.01448597695559319, 633, 671, 1123, 1286, 'AVD US Equity_EMN US Equity', $null),
    @(21, 0.06229374138390553, -0.0150605951347782, 0.03992897659002247, 659, 645, 951, 1458, 'AVD US Equity_FMC US Equity', $null),
    @(22, -0.02269950403380039, 0.002328271707850327, 0.02393637397455116, 661, 643, 1110, 1299, 'AVD US Equity_FUL US Equity', $null),
    @(23, -0.1878788166199025, 0.003058109627089634, 0.04880065627298124, 447, 579, 447, 380, 'AVD US Equity_GCP US Equity', $null),
    @(24, 0.0046436602133344, 0.00779186612164029, 0.0238366028293344, 623, 681, 1176, 1233, 'AVD US Equity_GRA US Equity', $null),
    @(25, -0.0813744242963016, 0.00288188710700199, 0.01215262093278261, 587, 717, 1122, 1287, 'AVD US Equity_IFF US Equity', $null),
    @(26, -0.1277347393047448, -0.005468922083822436, 0.1321521245122568, 584, 720, 1027, 1382, 'AVD US Equity_HUN US Equity', $null),
    @(27, -0.01881210072338735, -0.005349635170707301, 0.0295283505090411, 711, 593, 965, 1444, 'AVD US Equity_HXL US Equity', $null),
    @(28, 0.06396423374118054, -0.0007938215469115439, 0.07454582793185634, 502, 802, 1194, 1215, 'AVD US Equity_KOP US Equity', $null),
    @(29, 0.1622108468365961, -0.005231060483790961, 0.0122255863725198, 705, 599, 1107, 1302, 'AVD US Equity_LIN US Equity', $null),
    @(30, -0.01182962228869533, 0.004869429360998545, 0.01390432906132748, 583, 721, 1156, 1170, 'AVD US Equity_LYB US Equity', $null),
    @(31, -0.1435038081531846, 0.0025277971731541, 0.005677070886054307, 609, 695, 1175, 1234, 'AVD US Equity_MMM US Equity', $null),
    @(32, 0.005230150930100264, 0.004786082850642526, 0.03839924480669842, 598, 706, 1084, 1325, 'AVD US Equity_NZYMB DC Equity', $null),
    @(33, 0.1395712675510665, 0.005814618446992936, 0.004797838696924469, 583, 721, 1078, 1331, 'AVD US Equity_NEU US Equity', $null),
    @(34, -0.2167644904995195, -0.002012203588623751, 0.08912389304989515, 631, 673, 583, 636, 'AVD US Equity_OEC US Equity', $null),
    @(35, 0.1735576088160105, 0.004798596564365543, 0.01698331064492154, 664, 640, 1080, 1329, 'AVD US Equity_PPG US Equity', $null),
    @(36, -0.1247476419288321, -0.01007666768641635, 0.03350250553855097, 436, 868, 592, 1258, 'AVD US Equity_PRLB US Equity', $null),
    @(37, 0.2241411911575006, 0.0001509903550658276, 0.03994995511622357, 625, 679, 1127, 1282, 'AVD US Equity_RPM US Equity', $null),
    @(38, 0.1849806090513459, -0.008631376606841323, 0.005288584701666238, 682, 622, 966, 1443, 'AVD US Equity_SHW US Equity', $null),
    @(39, -0.4802265739462693, -0.02077412169385895, 0.05095635460197994, 597, 707, 485, 766, 'AVD US Equity_TSE US Equity', $null),
    @(40, -0.094382393282473, 0.0004903505226563354, 0.06670388508312587, 631, 554, 443, 543, 'AVD US Equity_UNVR US Equity', $null),
    @(41, 0.08226413576308222, -0.01057218849945696, 0.01513860508604028, 558, 746, 1006, 1403, 'AVD US Equity_WDFC US Equity', $null),
    @(42, -0.09555487622910819, -0.0008038894603350499, 0.0651705133059693, 403, 451, 320, 335, 'AVD US Equity_VVV US Equity', $null),
    @(43, -0.1541050413124814, -0.00154804113192053, 0.02923006328547544, 581, 723, 1140, 1269, 'AVD US Equity_WLK US Equity', $null),
    @(44, -0.2648502129603207, 0.004488406885562179, 0.1504868094317276, 475, 829, 772, 673, 'ADM US Equity_CINR US Equity', $null),
    @(45, -0.1522064647113875, -0.01302035881811103, 0.0440565866278962, 643, 661, 742, 1667, 'ADM US Equity_HXL US Equity', $null),
    @(46, -0.1237068352599696, -0.01370353751723619, 0.02831621764869633, 554, 750, 330, 812, 'AXTA US Equity_HXL US Equity', $null),
    @(47, -0.2054356968278096, 0.003718052595605981, 0.2105607291820612, 540, 764, 978, 467, 'CRDA LN Equity_CINR US Equity', $null),
    @(48, -0.03770052414040692, -0.002601026754238678, 0.05061990068986227, 759, 545, 1160, 1249, 'CRDA LN Equity_HXL US Equity', $null),
    @(49, -0.037612408018321, 0.007535285119400381, 0.05063558396495213, 670, 634, 800, 440, 'IMCD NA Equity_HXL US Equity', $null),
    @(50, -0.06612800862919, -0.0006072694588921124, 0.04456882871483179, 637, 667, 995, 1414, 'DCI US Equity_HXL US Equity', $null),
    @(51, -0.0809942996131332, 0.003351723959316821, 0.06260601494458734, 611, 693, 1142, 1267, 'DSM NA Equity_HXL US Equity', $null),
    @(52, 0.1982579926730657, 0.00628858209343508, 0.1041538314994525, 556, 748, 973, 1436, $null, 'CBT US Equity_POL US Equity'),
    @(53, -0.09696925015914404, -0.0002863403325866765, 0.1100657884578305, 671, 633, 1102, 1307, 'ECL US Equity_HXL US Equity', $null),
    @(54, -0.01826339385463838, -0.006167334564040727, 0.05938847966350897, 617, 687, 845, 1564, 'FUL US Equity_HXL US Equity', $null),
    @(55, 0.178578311466091, -0.0001275473744098976, 0.01654872645815975, 629, 675, 1393, 1016, $null, 'HXL US Equity_PPG US Equity'),
    @(56, 0.1216357816824767, 0.01442106087000616, 0.142621365248795, 682, 622, 1317, 1092, $null, 'HXL US Equity_RPM US Equity'),
    @(57, 0.2175673055895437, 0.01091176250744286, 0.02645582957958477, 656, 648, 1312, 1097, $null, 'HXL US Equity_SOLB BB Equity'),
    @(58, 0.25186347737743, 0.005770691267939387, 0.0716226857256511, 674, 630, 1322, 1087, $null, 'HXL US Equity_RPM US Equity'),
    @(59, 0.03942242597478973, 0.008131991361738589, 0.04486156542760934, 646, 658, 1546, 863, $null, 'HXL US Equity_SOLB BB Equity')
)

$styleSrc = $ws.Cells.Item(2, 1)

foreach ($row in $data) {
    $r = $row[0] + 2
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    if ($row[8] -ne $null) { $ws.Cells.Item($r, 9).Value = $row[8] } else { $ws.Cells.Item($r, 9).ClearContents() }
    if ($row[9] -ne $null) { $ws.Cells.Item($r, 10).Value = $row[9] } else { $ws.Cells.Item($r, 10).ClearContents() }
}

# Ensure column-A style (bold, centered, thin border) on the two newly-added rows
foreach ($rr in @(60, 61)) {
    $dst = $ws.Cells.Item($rr, 1)
    $dst.Borders.LineStyle = $styleSrc.Borders.LineStyle
    $dst.Font.Bold = $styleSrc.Font.Bold
    $dst.HorizontalAlignment = $styleSrc.HorizontalAlignment
    $dst.VerticalAlignment = $styleSrc.VerticalAlignment
}
